# Release Candidate: before Site check
#
# Updates a handful of data cells in the BLOCK_PATTERN table and trims the
# sheet's current selection back to the used range (A2:XFD26 instead of
# A2:XFD27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------
# First block: Duration[sec] value
$ws.Range("C4").Value = 1

# Second block: Duration[sec] value
$ws.Range("C15").Value = 1

# Second block: num waves from / to
$ws.Range("C17").Value = 0.5
$ws.Range("C18").Value = 1.5

# --- Selection --------------------------------------------------------------
# Shrink the selected range by one row (A2:XFD27 -> A2:XFD26) and keep the
# active cell at A2.
$ws.Range("A2:XFD26").Select()
